$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The survey rows for SQUIMAN (rows 137-151) were replaced by the rows that
# used to be further down the sheet (old rows 152-164, the SOLEVUL catch
# records). Deleting rows 137-151 with a shift-up accomplishes exactly that:
# the old row 152 becomes the new row 137, old row 164 becomes new row 149,
# and the sheet shrinks from 164 to 149 data rows.
$rng = $ws.Range("A137:Q151")
$rng.Delete()
